$d = $word.ActiveDocument

# 1. "Static objects (platforms, obstacles, etc.)" -> "More static objects (platforms, obstacles, etc.)"
#    (authored as two runs: "More s" + "tatic objects (platforms, obstacles, etc.)")
$d.Content.Find.Execute("Static objects (platforms, obstacles, etc.)", $true, $false, $false, $false, $false, $true, 1, $false, "More static objects (platforms, obstacles, etc.)", 2)

$rng = $d.Content
$rng.Find.Execute("More static objects (platforms, obstacles, etc.)")
$splitPoint = $rng.Start

$tempBookmark = $d.Bookmarks.Add("tempSplitMark", $d.Range($splitPoint, $splitPoint + 6))
$d.Bookmarks("tempSplitMark").Delete()

# 2. Add two new bullet items after "Game over (run out of lives)"
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara1 = $d.Paragraphs.Last
$newPara1.Range.Text = "Make level longer in general"

$newPara1.Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs.Last
$newPara2.Range.Text = "Maybe high score list?"
